# Scheduled market-data refresh: updates currentAveragePrice / Leve profit
# columns (H-N) across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 1000
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0

# Row 129
$ws.Range("H129").Value = 1130.2683
$ws.Range("J129").Value = 1261.7428
$ws.Range("L129").Value = 3785.2284
$ws.Range("N129").Value = -13785.2284

# Row 132
$ws.Range("H132").Value = 3191.5715
$ws.Range("I132").Value = 3228.7144
$ws.Range("J132").Value = 3043
$ws.Range("K132").Value = 9686.143199999999
$ws.Range("L132").Value = 9129
$ws.Range("M132").Value = -7156.143199999999
$ws.Range("N132").Value = -14189

# Row 137
$ws.Range("H137").Value = 1761.0435
$ws.Range("I137").Value = 1715.375
$ws.Range("K137").Value = 5146.125
$ws.Range("M137").Value = -2596.125

# Row 138
$ws.Range("H138").Value = 50004784
$ws.Range("I138").Value = 200005000
$ws.Range("J138").Value = 4710.6665
$ws.Range("K138").Value = 600015000
$ws.Range("L138").Value = 14131.9995
$ws.Range("M138").Value = -600009860
$ws.Range("N138").Value = -24411.9995

# Row 141
$ws.Range("H141").Value = 2184.1853
$ws.Range("I141").Value = 1468.7
$ws.Range("K141").Value = 4406.1
$ws.Range("M141").Value = 773.8999999999996

$ws = $wb.Worksheets.Item("ARM")
# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21538

# Row 74
$ws.Range("H74").Value = 142862640
$ws.Range("I74").Value = 500006240
$ws.Range("J74").Value = 5199.8
$ws.Range("K74").Value = 500006240
$ws.Range("L74").Value = 5199.8
$ws.Range("M74").Value = -500005366
$ws.Range("N74").Value = -6947.8

# Row 77
$ws.Range("H77").Value = 142862640
$ws.Range("I77").Value = 500006240
$ws.Range("J77").Value = 5199.8
$ws.Range("K77").Value = 2500031200
$ws.Range("L77").Value = 25999
$ws.Range("M77").Value = -2500026832
$ws.Range("N77").Value = -34735

# Row 102
$ws.Range("H102").Value = 990.6923
$ws.Range("I102").Value = 990.6923
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 990.6923
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 631.3077

# Row 132
$ws.Range("H132").Value = 14495.125
$ws.Range("I132").Value = 1813.4642
$ws.Range("J132").Value = 44085.668
$ws.Range("K132").Value = 5440.392599999999
$ws.Range("L132").Value = 132257.004
$ws.Range("M132").Value = -2910.392599999999
$ws.Range("N132").Value = -137317.004

$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

# Row 134
$ws.Range("H134").Value = 1875
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 609.8889
$ws.Range("I26").Value = 518
$ws.Range("J26").Value = 724.75
$ws.Range("K26").Value = 1554
$ws.Range("L26").Value = 2174.25
$ws.Range("M26").Value = -1266
$ws.Range("N26").Value = -2750.25

# Row 80
$ws.Range("H80").Value = 6426.316
$ws.Range("J80").Value = 6616.722
$ws.Range("L80").Value = 19850.166
$ws.Range("N80").Value = -21722.166

# Row 83
$ws.Range("H83").Value = 6426.316
$ws.Range("J83").Value = 6616.722
$ws.Range("L83").Value = 59550.498
$ws.Range("N83").Value = -68910.49799999999

# Row 131
$ws.Range("H131").Value = 706.42
$ws.Range("J131").Value = 730.35486
$ws.Range("L131").Value = 2191.06458
$ws.Range("N131").Value = -12271.06458

$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 20005600
$ws.Range("J52").Value = 20005600
$ws.Range("L52").Value = 20005600
$ws.Range("N52").Value = -20006118

# Row 53
$ws.Range("H53").Value = 5039
$ws.Range("I53").Value = 5039
$ws.Range("K53").Value = 5039
$ws.Range("M53").Value = -4408

# Row 102
$ws.Range("H102").Value = 38465644
$ws.Range("I102").Value = 45458348
$ws.Range("K102").Value = 45458348
$ws.Range("M102").Value = -45456726

# Row 122
$ws.Range("H122").Value = 78432930
$ws.Range("I122").Value = 23811164
$ws.Range("K122").Value = 71433492
$ws.Range("M122").Value = -71431042

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 11972.143
$ws.Range("I7").Value = 5525
$ws.Range("J7").Value = 20568.334
$ws.Range("K7").Value = 5525
$ws.Range("L7").Value = 20568.334
$ws.Range("M7").Value = -5413
$ws.Range("N7").Value = -20792.334

# Row 16
$ws.Range("H16").Value = 639.2727
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340

# Row 22
$ws.Range("H22").Value = 4212.375
$ws.Range("I22").Value = 5300
$ws.Range("J22").Value = 3849.8333
$ws.Range("K22").Value = 5300
$ws.Range("L22").Value = 3849.8333
$ws.Range("M22").Value = -5005
$ws.Range("N22").Value = -4439.8333

# Row 27
$ws.Range("H27").Value = 4212.375
$ws.Range("I27").Value = 5300
$ws.Range("J27").Value = 3849.8333
$ws.Range("K27").Value = 5300
$ws.Range("L27").Value = 3849.8333
$ws.Range("M27").Value = -5193
$ws.Range("N27").Value = -4063.8333

# Row 46
$ws.Range("H46").Value = 641.56525
$ws.Range("I46").Value = 611.75
$ws.Range("J46").Value = 674.0909
$ws.Range("K46").Value = 611.75
$ws.Range("L46").Value = 674.0909
$ws.Range("M46").Value = -423.75
$ws.Range("N46").Value = -1050.0909

# Row 122
$ws.Range("H122").Value = 1228527.2
$ws.Range("I122").Value = 1785067.1
$ws.Range("J122").Value = 4139.8
$ws.Range("K122").Value = 5355201.300000001
$ws.Range("L122").Value = 12419.4
$ws.Range("M122").Value = -5352751.300000001
$ws.Range("N122").Value = -17319.4

# Row 126
$ws.Range("H126").Value = 11972.143
$ws.Range("I126").Value = 5525
$ws.Range("J126").Value = 20568.334
$ws.Range("K126").Value = 16575
$ws.Range("L126").Value = 61705.00199999999
$ws.Range("M126").Value = -14105
$ws.Range("N126").Value = -66645.00199999999

# Row 132
$ws.Range("H132").Value = 503763.72
$ws.Range("I132").Value = 710081.5
$ws.Range("J132").Value = 2706.1428
$ws.Range("K132").Value = 2130244.5
$ws.Range("L132").Value = 8118.428400000001
$ws.Range("M132").Value = -2127714.5
$ws.Range("N132").Value = -13178.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1319.8
$ws.Range("I122").Value = 1349.8334
$ws.Range("J122").Value = 1199.6666
$ws.Range("K122").Value = 4049.5002
$ws.Range("L122").Value = 3598.9998
$ws.Range("M122").Value = -1599.5002
$ws.Range("N122").Value = -8498.9998

# Row 128
$ws.Range("H128").Value = 38750
$ws.Range("J128").Value = 38750
$ws.Range("L128").Value = 38750
$ws.Range("N128").Value = -48710
